# Apply the "new version with timestamp" update to the day-sale report.
# Rows 9-13 (products list) shift: the "BRUFEN 400MG 30 TAB" row is replaced,
# the rows below move up by one, and a new "OTAL EAR DROPS 5 ML" row appears
# at row 13. The grand-total cell and the printed timestamp are refreshed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Columns such as P are formatted with a numeric NumberFormat (e.g. "0.00")
    # even though the underlying cell actually stores a text string in this
    # workbook. Assigning a numeric-looking string directly would make Excel
    # silently coerce it to a real number, which changes the cell's type in
    # the saved XML. Temporarily switching to a text format preserves the
    # original "numeric-looking text" representation.
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# Row 9 : BRUFEN 400MG 30 TAB  ->  EREC 100MG 12 F.C. TABLETS
$ws.Range("C9").Value = "EREC 100MG 12 F.C. TABLETS"
$ws.Range("H9").Value = "2:7"
$ws.Range("N9").Value = "144.00"
Set-TextValue $ws.Range("P9") "11.5200"

# Row 10 : EREC 100MG 12 F.C. TABLETS  ->  FLECTOR 50MG 30 CAPS
$ws.Range("C10").Value = "FLECTOR 50MG 30 CAPS"
$ws.Range("H10").Value = "1:1"
$ws.Range("N10").Value = "87.00"
Set-TextValue $ws.Range("P10") "28.7100"

# Row 11 : FLECTOR 50MG 30 CAPS  ->  ICANDRA 50 MG 30 TAB.
$ws.Range("C11").Value = "ICANDRA 50 MG 30 TAB."
$ws.Range("H11").Value = "0:2"
$ws.Range("N11").Value = "118.50"
Set-TextValue $ws.Range("P11") "78.2100"
$ws.Range("Q11").Value = "0:2"

# Row 12 : ICANDRA 50 MG 30 TAB.  ->  KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF.
$ws.Range("C12").Value = "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF."
$ws.Range("H12").Value = "1:1"
$ws.Range("N12").Value = "60.00"
Set-TextValue $ws.Range("P12") "60.0000"
$ws.Range("Q12").Value = "1:0"

# Row 13 : KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF.  ->  OTAL EAR DROPS 5 ML (new row)
$ws.Range("C13").Value = "OTAL EAR DROPS 5 ML"
$ws.Range("H13").Value = "5:0"
$ws.Range("N13").Value = "19.00"
Set-TextValue $ws.Range("P13") "19.0000"

# Grand total (sum of price column) decreases because 25.74 was replaced by 19.00
$ws.Range("P23").Value = 476.24

# Refresh the printed report timestamp shown in the footer
$ws.Range("A24").Value = "Friday, 29 August, 2025 5:07 PM"
